# TimeLog_ConnorPeper.xlsx update
# - Append a period to the end of the "Task 22 ... Task 17" activity text (row 13, column F)
# - Fill in row 14 (week 8) with Hours = 2 and Activities = "Task 16: Complete."
# - Move the active selection from E14 to E15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the trailing punctuation on the prior week's activity note (F13)
$ws.Range("F13").Value = "Task 22: Got some basic error checking working. Can catch the exceptions though right now they're just ignored. Task 13: Taught Shammi how to make a Stored Procedure, Throw Errors, Do IF statements, and store variables in tSQL. Task 17: Created the SPROC to insert a reply into the database."

# Record hours worked and activity note for week 8 (row 14)
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = "Task 16: Complete."

# Move selection to E15, matching the saved workbook view
$ws.Range("E15").Select()
